# Regenerate merged AHB files
# 1) Rename the "_old" / "_new" header-name suffixes to "_FV2404" / "_FV2410"
# 2) Turn the data range into a real table (Table1)
# 3) Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $baseNames[$i] + "_FV2404"
    $ws.Cells.Item(1, $i + 12).Value2 = $baseNames[$i] + "_FV2410"
}

# Turn A1:U72 into a proper table with an autofilter
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U72"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false

# Freeze the header row (row 1) in the sheet view
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
